$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hide columns A:B (the two frozen, translated "comment" columns)
$ws.Columns("A:B").Hidden = $true

# Move the selection in the (frozen-pane) sheet view to G10
[void]$ws.Range("G10").Select()
